$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 8360.799999999999
$ws.Range("I28").Value = 7951.625
$ws.Range("K28").Value = 7951.625
$ws.Range("M28").Value = -7466.625

$ws.Range("H41").Value = 246.66667
$ws.Range("I41").Value = 246.66667
$ws.Range("K41").Value = 246.66667
$ws.Range("M41").Value = 193.33333

$ws.Range("H51").Value = 9999.799999999999
$ws.Range("I51").Value = 9999
$ws.Range("K51").Value = 9999
$ws.Range("M51").Value = -9515

$ws.Range("H70").Value = 155846.95
$ws.Range("I70").Value = 223361.08
$ws.Range("J70").Value = 9566.333000000001
$ws.Range("K70").Value = 670083.24
$ws.Range("L70").Value = 28698.999
$ws.Range("M70").Value = -669813.24
$ws.Range("N70").Value = -29238.999

$ws.Range("H73").Value = 155846.95
$ws.Range("I73").Value = 223361.08
$ws.Range("J73").Value = 9566.333000000001
$ws.Range("K73").Value = 670083.24
$ws.Range("L73").Value = 28698.999
$ws.Range("M73").Value = -669147.24
$ws.Range("N73").Value = -30570.999

$ws.Range("H113").Value = 6036.6924
$ws.Range("I113").Value = 5634.727
$ws.Range("K113").Value = 5634.727
$ws.Range("M113").Value = -2380.727

$ws.Range("H132").Value = 502167.34
$ws.Range("I132").Value = 1961.6471
$ws.Range("K132").Value = 5884.9413
$ws.Range("M132").Value = -3354.9413

$ws.Range("H138").Value = 3873.3333
$ws.Range("I138").Value = 2168.3157
$ws.Range("J138").Value = 4725.8423
$ws.Range("K138").Value = 6504.9471
$ws.Range("L138").Value = 14177.5269
$ws.Range("M138").Value = -1364.9471
$ws.Range("N138").Value = -24457.5269

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 4771.4
$ws.Range("J14").Value = 5914.25
$ws.Range("L14").Value = 5914.25
$ws.Range("N14").Value = -6264.25

$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H32").Value = 8787.462
$ws.Range("I32").Value = 4597.727
$ws.Range("J32").Value = 31831
$ws.Range("K32").Value = 4597.727
$ws.Range("L32").Value = 31831
$ws.Range("M32").Value = -4310.727
$ws.Range("N32").Value = -32405

$ws.Range("H33").Value = 3296.5715
$ws.Range("J33").Value = 9500
$ws.Range("L33").Value = 9500
$ws.Range("N33").Value = -10158

$ws.Range("H61").Value = 3973.5557
$ws.Range("I61").Value = 3705.05
$ws.Range("K61").Value = 3705.05
$ws.Range("M61").Value = -3493.05

$ws.Range("H88").Value = 2231.1177
$ws.Range("I88").Value = 1557.6666
$ws.Range("J88").Value = 2598.4546
$ws.Range("K88").Value = 1557.6666
$ws.Range("L88").Value = 2598.4546
$ws.Range("M88").Value = -1151.6666
$ws.Range("N88").Value = -3410.4546

$ws.Range("H91").Value = 2231.1177
$ws.Range("I91").Value = 1557.6666
$ws.Range("J91").Value = 2598.4546
$ws.Range("K91").Value = 1557.6666
$ws.Range("L91").Value = 2598.4546
$ws.Range("M91").Value = -153.6666
$ws.Range("N91").Value = -5406.4546

$ws.Range("H110").Value = 1620.4736
$ws.Range("I110").Value = 1010
$ws.Range("K110").Value = 1010
$ws.Range("M110").Value = 1035

$ws.Range("H132").Value = 868.75
$ws.Range("J132").Value = 5119.5
$ws.Range("L132").Value = 15358.5
$ws.Range("N132").Value = -20418.5

$ws.Range("H136").Value = 3973.5557
$ws.Range("I136").Value = 3705.05
$ws.Range("K136").Value = 11115.15
$ws.Range("M136").Value = -8565.150000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2399.7
$ws.Range("I94").Value = 2001
$ws.Range("J94").Value = 2665.5
$ws.Range("K94").Value = 2001
$ws.Range("L94").Value = 2665.5
$ws.Range("M94").Value = -1550
$ws.Range("N94").Value = -3567.5

$ws.Range("H134").Value = 2487.5789
$ws.Range("I134").Value = 2140.9375
$ws.Range("J134").Value = 4336.3335
$ws.Range("K134").Value = 6422.8125
$ws.Range("L134").Value = 13009.0005
$ws.Range("M134").Value = -3887.8125
$ws.Range("N134").Value = -18079.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4277.132
$ws.Range("I31").Value = 2170.75
$ws.Range("J31").Value = 6636.28
$ws.Range("K31").Value = 2170.75
$ws.Range("L31").Value = 6636.28
$ws.Range("M31").Value = -1875.75
$ws.Range("N31").Value = -7226.28

$ws.Range("H34").Value = 4277.132
$ws.Range("I34").Value = 2170.75
$ws.Range("J34").Value = 6636.28
$ws.Range("K34").Value = 2170.75
$ws.Range("L34").Value = 6636.28
$ws.Range("M34").Value = -1968.75
$ws.Range("N34").Value = -7040.28

$ws.Range("H52").Value = 73966.336
$ws.Range("J52").Value = 73966.336
$ws.Range("L52").Value = 73966.336
$ws.Range("N52").Value = -74554.336

$ws.Range("H86").Value = 3539.7273
$ws.Range("I86").Value = 3210.5
$ws.Range("J86").Value = 4417.6665
$ws.Range("K86").Value = 3210.5
$ws.Range("L86").Value = 4417.6665
$ws.Range("M86").Value = -2087.5
$ws.Range("N86").Value = -6663.6665

$ws.Range("H89").Value = 3539.7273
$ws.Range("I89").Value = 3210.5
$ws.Range("J89").Value = 4417.6665
$ws.Range("K89").Value = 16052.5
$ws.Range("L89").Value = 22088.3325
$ws.Range("M89").Value = -10436.5
$ws.Range("N89").Value = -33320.3325

$ws.Range("H99").Value = 2209.238
$ws.Range("I99").Value = 2475.5334
$ws.Range("K99").Value = 2475.5334
$ws.Range("M99").Value = -977.5333999999998

$ws.Range("H107").Value = 530.8946999999999
$ws.Range("I107").Value = 518.8889
$ws.Range("J107").Value = 541.7
$ws.Range("K107").Value = 518.8889
$ws.Range("L107").Value = 541.7
$ws.Range("M107").Value = 1401.1111
$ws.Range("N107").Value = -4381.7

$ws.Range("H122").Value = 1984.25
$ws.Range("I122").Value = 1984.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5952.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3502.75
$ws.Range("N122").ClearContents()

$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 30000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -34920

$ws.Range("H126").Value = 2209.238
$ws.Range("I126").Value = 2475.5334
$ws.Range("K126").Value = 7426.600199999999
$ws.Range("M126").Value = -4956.600199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 248.33333
$ws.Range("J46").Value = 248.33333
$ws.Range("L46").Value = 744.99999
$ws.Range("N46").Value = -926.99999

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1144.35
$ws.Range("I97").Value = 1046.7894
$ws.Range("K97").Value = 1046.7894
$ws.Range("M97").Value = -550.7893999999999

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 4983.5884
$ws.Range("I126").Value = 4729.625
$ws.Range("J126").Value = 5209.3335
$ws.Range("K126").Value = 14188.875
$ws.Range("L126").Value = 15628.0005
$ws.Range("M126").Value = -11718.875
$ws.Range("N126").Value = -20568.0005

$ws.Range("H132").Value = 3679.1667
$ws.Range("I132").Value = 3195.3333
$ws.Range("K132").Value = 9585.999899999999
$ws.Range("M132").Value = -7055.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 45000
$ws.Range("J75").Value = 45000
$ws.Range("L75").Value = 45000
$ws.Range("N75").Value = -46872

$ws.Range("H78").Value = 45000
$ws.Range("J78").Value = 45000
$ws.Range("L78").Value = 135000
$ws.Range("N78").Value = -144360

$ws.Range("H82").Value = 1858.8
$ws.Range("J82").Value = 2698.2
$ws.Range("L82").Value = 2698.2
$ws.Range("N82").Value = -3420.2

$ws.Range("H85").Value = 1858.8
$ws.Range("J85").Value = 2698.2
$ws.Range("L85").Value = 2698.2
$ws.Range("N85").Value = -5194.2

$ws.Range("H122").Value = 5015
$ws.Range("I122").Value = 5400
$ws.Range("J122").Value = 4982.9165
$ws.Range("K122").Value = 16200
$ws.Range("L122").Value = 14948.7495
$ws.Range("M122").Value = -13750
$ws.Range("N122").Value = -19848.7495

$ws.Range("H132").Value = 5267.25
$ws.Range("J132").Value = 6069
$ws.Range("L132").Value = 18207
$ws.Range("N132").Value = -23267

$ws.Range("H136").Value = 3355.1904
$ws.Range("I136").Value = 3155.75
$ws.Range("J136").Value = 3993.4
$ws.Range("K136").Value = 9467.25
$ws.Range("L136").Value = 11980.2
$ws.Range("M136").Value = -6917.25
$ws.Range("N136").Value = -17080.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 2200
$ws.Range("I23").Value = 2200
$ws.Range("K23").Value = 2200
$ws.Range("M23").Value = -1971

$ws.Range("H100").Value = 938.3333
$ws.Range("I100").Value = 785
$ws.Range("J100").Value = 1475
$ws.Range("K100").Value = 1570
$ws.Range("L100").Value = 2950
$ws.Range("M100").Value = -1029
$ws.Range("N100").Value = -4032

$ws.Range("H113").Value = 1128.6666
$ws.Range("I113").Value = 807.6667
$ws.Range("K113").Value = 2423.0001
$ws.Range("M113").Value = -253.0001000000002

$ws.Range("H122").Value = 3416.5806
$ws.Range("I122").Value = 2943.0435
$ws.Range("K122").Value = 8829.130500000001
$ws.Range("M122").Value = -6379.130500000001

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H132").Value = 4330.6
$ws.Range("I132").Value = 4330.6
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12991.8
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10461.8
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 1998.7693
$ws.Range("I136").Value = 1807.3914
$ws.Range("J136").Value = 3466
$ws.Range("K136").Value = 5422.174199999999
$ws.Range("L136").Value = 10398
$ws.Range("M136").Value = -2872.174199999999
$ws.Range("N136").Value = -15498
